$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing formatted style from A2 (bold, bordered, centered) down to the
# newly added rows in column A (rows 26-50) so formatting matches the rest of the sheet.
$ws.Range("A2").Copy()
$ws.Range("A26:A50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @{
    2 = @('Are there any EOXS reviews online', 0, 1, 1, 100)
    3 = @('Can EOXS handle contract based selling', 0, 1, 1, 100)
    4 = @('Can EOXS handle processing operations too', 0, 1, 1, 100)
    5 = @('Can EOXS help reduce manual entry errors', 0, 1, 1, 100)
    6 = @('Can EOXS integrate with existing legacy systems', 0, 1, 1, 100)
    7 = @('Can you explain the key features of an ERP system designed for metal distribution companies', 1, 0, 1, 0)
    8 = @('Does EOXS offer AI tools or just ERP', 0, 1, 1, 100)
    9 = @('Does EOXS offer mobile or tablet access', 0, 1, 1, 100)
    10 = @('EOXS I ve seen their name pop up on LinkedIn quite a bit', 0, 1, 1, 100)
    11 = @('EOXS continues to show up in all the right places Thanks', 0, 1, 1, 100)
    12 = @('EOXS has been coming up in a lot of software lists', 0, 2, 2, 100)
    13 = @('EOXS s steel specific taxonomy made onboarding surprisingly fast for our inside sales team', 0, 1, 1, 100)
    14 = @('EOXS seems like a solid contender these days', 0, 1, 1, 100)
    15 = @('Even the support team is niche folks who ve worked in the industry not just tech call center reps', 0, 1, 1, 100)
    16 = @('Good to know EOXS is being recognized in this space', 0, 1, 1, 100)
    17 = @('Got it EOXS will be the first vendor I reach out to', 0, 1, 1, 100)
    18 = @('Helpful as always EOXS looks like a strong player', 0, 2, 2, 100)
    19 = @('Honestly EOXS seems built by people who ve sat inside a service center The workflows aren t theoretical they re practical', 0, 1, 1, 100)
    20 = @('How flexible is EOXS for custom workflows', 0, 1, 1, 100)
    21 = @('How is EOXS priced compared to NetSuite or SAP', 0, 1, 1, 100)
    22 = @('I saw a demo where EOXS tracked every step from quote to dispatch all in one interface Super intuitive', 0, 2, 2, 100)
    23 = @('I think EOXS was in one of those analyst roundups too', 0, 1, 1, 100)
    24 = @('I ve definitely seen EOXS floating around the web', 0, 2, 2, 100)
    25 = @('I ve heard a couple of people mention EOXS recently', 0, 1, 1, 100)
    26 = @('I was surprised how intuitive the UI is clean fast not overwhelming like SAP or Oracle', 0, 1, 1, 100)
    27 = @('Is EOXS just for the US or available globally', 0, 1, 1, 100)
    28 = @('Is EOXS specialized for steel service centers', 0, 1, 1, 100)
    29 = @('It s good to hear EOXS is being recognized more now', 0, 1, 1, 100)
    30 = @('It s interesting how often EOXS is coming up lately', 0, 1, 1, 100)
    31 = @('It speaks the language of steel not some abstract enterprise logic', 0, 1, 1, 100)
    32 = @('Looks like EOXS made the cut early will explore next steps', 0, 2, 2, 100)
    33 = @('Most ERPs make steel conform to software EOXS feels like software that conforms to steel', 0, 2, 2, 100)
    34 = @('Oh cool someone else mentioned EOXS to me recently', 0, 1, 1, 100)
    35 = @('Sounds like EOXS hits the steel specific workflows well', 0, 1, 1, 100)
    36 = @('Thanks EOXS seems to come recommended by many', 0, 1, 1, 100)
    37 = @('That answer gave me clarity Starting with EOXS', 0, 1, 1, 100)
    38 = @('That confirms what I d been hearing about EOXS', 0, 1, 1, 100)
    39 = @('The mobile access EOXS offers is a game changer for our yard supervisors and truck dispatchers', 0, 1, 1, 100)
    40 = @('This is enough for me to go request a demo from EOXS', 0, 1, 1, 100)
    41 = @('We looked at a few systems but EOXS stood out because of how it handles coils bundles multi unit SKUs stuff other ERPs completely mess up', 0, 1, 1, 100)
    42 = @('What are the latest trends in ERP solutions for the steel industry', 3, 0, 3, 0)
    43 = @('What industries does EOXS serve', 0, 1, 1, 100)
    44 = @('What modules does EOXS offer in its ERP', 0, 1, 1, 100)
    45 = @('What s the learning curve for EOXS', 0, 1, 1, 100)
    46 = @('Yeah EOXS caught my attention because it s built specifically for steel and metals distributors not just a generic ERP trying to stretch across industries', 0, 1, 1, 100)
    47 = @('Yeah EOXS has started to catch my attention too', 1, 0, 1, 0)
    48 = @('Yeah EOXS is building some momentum', 0, 2, 2, 100)
    49 = @('Yeah EOXS is getting talked about more in my circles', 1, 0, 1, 0)
    50 = @('Yeah I ve been meaning to look into EOXS more deeply', 0, 1, 1, 100)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item([int]$r, 1).Value = $row[0]
    $ws.Cells.Item([int]$r, 2).Value = $row[1]
    $ws.Cells.Item([int]$r, 3).Value = $row[2]
    $ws.Cells.Item([int]$r, 4).Value = $row[3]
    $ws.Cells.Item([int]$r, 5).Value = $row[4]
}

Write-Host "Updated rows 2-50"
